# Main.xlsx / "Rules" sheet: cell B11 used to hold the text "R40"
# (a rule name) but should now hold the text "1" instead. The cell
# must keep its original style/number-format and stay a *text* cell
# (t="s", referencing a new shared string "1") rather than turn into
# the number 1 - a plain Range.Value/Formula assignment of "1" would
# be auto-coerced by Excel into a numeric value, which is not what the
# edit calls for.
#
# To get a genuinely text-typed "1" without touching B11's
# NumberFormat (any NumberFormat/quote-prefix write on the cell itself
# mints a brand-new style record, which would change its style index),
# stage the text in a scratch cell via the TEXT() formula - whose
# result is always a string - then copy/paste only the *value* into
# B11. PasteSpecial(xlPasteValues) carries the text-typed result across
# without bringing any formatting with it, so B11's own style is left
# completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

$helper = $ws.Range("Z1")
$helper.Formula = '=TEXT(1,"0")'
$helper.Copy()

$target.PasteSpecial(-4163)  # xlPasteValues

$helper.Clear()
